$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$row = 8

$ws.Cells.Item($row, 1).Value = "15c484ef-ab3d-48d1-a1b8-6c0ba8fa2aaa"
$ws.Cells.Item($row, 2).Value = "gtrouy[ping test by week"
$ws.Cells.Item($row, 4).Value = "To Do"

# DueDate / Created look like dates/timestamps - force them to be stored
# as literal text (matching the rest of the column) instead of letting
# Excel auto-convert them into date serial numbers.
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "2025-05-26"
$ws.Cells.Item($row, 5).ClearFormats()

$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "2025-05-24 10:49:03"
$ws.Cells.Item($row, 6).ClearFormats()
